# Fruta / hortaliza, semanal
# A new weekly price record is inserted at the top of the data block
# (row 148), pushing the existing rows 148-161 down to 149-162 and
# extending the used range from A1:T161 to A1:T162.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 148; this shifts rows 148:161 down
# to 149:162 (matches the diff: the old row 148 record now lives at 149,
# ..., the old row 160 record now lives at 161, and the old row 161
# record now lives at 162).
$ws.Rows.Item(148).Insert()

# Populate the newly-inserted row 148 with the new record.
$ws.Cells.Item(148, 1).Value = 6
$ws.Cells.Item(148, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(148, 3).Value = "Metropolitana"
$ws.Cells.Item(148, 4).Value = 44476
$ws.Cells.Item(148, 5).Value = 13
$ws.Cells.Item(148, 6).Value = "Fruta"
$ws.Cells.Item(148, 7).Value = 100101
$ws.Cells.Item(148, 8).Value = "Berries"
$ws.Cells.Item(148, 9).Value = 100101001
$ws.Cells.Item(148, 10).Value = "Arándano (blue)"
$ws.Cells.Item(148, 11).Value = "Sin especificar"
$ws.Cells.Item(148, 12).Value = "Especial"
$ws.Cells.Item(148, 13).Value = 750
$ws.Cells.Item(148, 14).Value = 14000
$ws.Cells.Item(148, 15).Value = 14000
$ws.Cells.Item(148, 16).Value = 14000
$ws.Cells.Item(148, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(148, 18).Value = "Perú"
$ws.Cells.Item(148, 19).Value = 7000
$ws.Cells.Item(148, 20).Value = 2
